$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 76927896
$ws.Range("I62").Value = 250003000
$ws.Range("J62").Value = 5630.222
$ws.Range("K62").Value = 250003000
$ws.Range("L62").Value = 5630.222
$ws.Range("M62").Value = -250002376
$ws.Range("N62").Value = -6878.222

$ws.Range("H65").Value = 76927896
$ws.Range("I65").Value = 250003000
$ws.Range("J65").Value = 5630.222
$ws.Range("K65").Value = 1250015000
$ws.Range("L65").Value = 28151.11
$ws.Range("M65").Value = -1250011880
$ws.Range("N65").Value = -34391.11

$ws.Range("H132").Value = 6280.552
$ws.Range("I132").Value = 1750.2609
$ws.Range("K132").Value = 5250.7827
$ws.Range("M132").Value = -2720.7827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 18183258
$ws.Range("I45").Value = 45455144
$ws.Range("K45").Value = 45455144
$ws.Range("M45").Value = -45454767

$ws.Range("H74").Value = 34888.56
$ws.Range("I74").Value = 45037.914
$ws.Range("J74").Value = 13667.182
$ws.Range("K74").Value = 45037.914
$ws.Range("L74").Value = 13667.182
$ws.Range("M74").Value = -44163.914
$ws.Range("N74").Value = -15415.182

$ws.Range("H77").Value = 34888.56
$ws.Range("I77").Value = 45037.914
$ws.Range("J77").Value = 13667.182
$ws.Range("K77").Value = 225189.57
$ws.Range("L77").Value = 68335.91
$ws.Range("M77").Value = -220821.57
$ws.Range("N77").Value = -77071.91

$ws.Range("H122").Value = 1515.619
$ws.Range("I122").Value = 1448.3077
$ws.Range("J122").Value = 1625
$ws.Range("K122").Value = 4344.9231
$ws.Range("L122").Value = 4875
$ws.Range("M122").Value = -1894.9231
$ws.Range("N122").Value = -9775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 955666.5
$ws.Range("I134").Value = 1145601.6
$ws.Range("J134").Value = 5990.5713
$ws.Range("K134").Value = 3436804.8
$ws.Range("L134").Value = 17971.7139
$ws.Range("M134").Value = -3434269.8
$ws.Range("N134").Value = -23041.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = ""
$ws.Range("N41").Value = -10856

$ws.Range("H134").Value = 4116
$ws.Range("I134").Value = 1085.4286
$ws.Range("J134").Value = 7651.6665
$ws.Range("K134").Value = 3256.2858
$ws.Range("L134").Value = 22954.9995
$ws.Range("M134").Value = -721.2857999999997
$ws.Range("N134").Value = -28024.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 12992
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 12992
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 12992
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = -13218

$ws.Range("H7").Value = 2433.3333
$ws.Range("I7").Value = 1650
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1650
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1538
$ws.Range("N7").Value = -4224

$ws.Range("H8").Value = 2433.3333
$ws.Range("I8").Value = 1650
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 1650
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = -1511
$ws.Range("N8").Value = -4278

$ws.Range("H9").Value = 15789
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 15789
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 15789
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = -16129

$ws.Range("H10").Value = 10526
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 10526
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 10526
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -10864

$ws.Range("H16").Value = 12992
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 12992
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 12992
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = -13492

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""

$ws.Range("H21").Value = 6000
$ws.Range("J21").Value = 6000
$ws.Range("L21").Value = 6000
$ws.Range("N21").Value = -6346

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

$ws.Range("H30").Value = 6000
$ws.Range("J30").Value = 6000
$ws.Range("L30").Value = 6000
$ws.Range("N30").Value = -6210

$ws.Range("H33").Value = 6673
$ws.Range("J33").Value = 6673
$ws.Range("L33").Value = 6673
$ws.Range("N33").Value = -7177

$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""

$ws.Range("H40").Value = 7250
$ws.Range("J40").Value = 7250
$ws.Range("L40").Value = 7250
$ws.Range("N40").Value = -7552

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""

$ws.Range("H46").Value = 25023
$ws.Range("J46").Value = 25023
$ws.Range("L46").Value = 25023
$ws.Range("N46").Value = -25335

$ws.Range("H57").Value = 19025
$ws.Range("J57").Value = 19025
$ws.Range("L57").Value = 19025
$ws.Range("N57").Value = -20665

$ws.Range("H122").Value = 1055
$ws.Range("I122").Value = 909.4
$ws.Range("J122").Value = 1601
$ws.Range("K122").Value = 2728.2
$ws.Range("L122").Value = 4803
$ws.Range("M122").Value = -278.1999999999998
$ws.Range("N122").Value = -9703

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 11025
$ws.Range("J50").Value = 11025
$ws.Range("L50").Value = 11025
$ws.Range("N50").Value = -12299

$ws.Range("H51").Value = 11332.667
$ws.Range("J51").Value = 11332.667
$ws.Range("L51").Value = 11332.667
$ws.Range("N51").Value = -12288.667

$ws.Range("H53").Value = 11985.571
$ws.Range("J53").Value = 13750
$ws.Range("L53").Value = 13750
$ws.Range("N53").Value = -14786

$ws.Range("H58").Value = 24998
$ws.Range("I58").Value = 24998
$ws.Range("K58").Value = 24998
$ws.Range("M58").Value = -24738

$ws.Range("H61").Value = 1980.6757
$ws.Range("I61").Value = 1794.2222
$ws.Range("J61").Value = 2484.1
$ws.Range("K61").Value = 1794.2222
$ws.Range("L61").Value = 2484.1
$ws.Range("M61").Value = -1592.2222
$ws.Range("N61").Value = -2888.1

$ws.Range("H100").Value = 1432.7778
$ws.Range("I100").Value = 1360
$ws.Range("J100").Value = 1460.7693
$ws.Range("K100").Value = 1360
$ws.Range("L100").Value = 1460.7693
$ws.Range("M100").Value = -819
$ws.Range("N100").Value = -2542.7693

$ws.Range("H113").Value = 1980.6757
$ws.Range("I113").Value = 1794.2222
$ws.Range("J113").Value = 2484.1
$ws.Range("K113").Value = 1794.2222
$ws.Range("L113").Value = 2484.1
$ws.Range("M113").Value = 375.7778000000001
$ws.Range("N113").Value = -6824.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1126.0952
$ws.Range("I126").Value = 1068
$ws.Range("J126").Value = 1220.5
$ws.Range("K126").Value = 3204
$ws.Range("L126").Value = 3661.5
$ws.Range("M126").Value = -734
$ws.Range("N126").Value = -8601.5
